$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B7:C8 raw values (Junho / Julho collections were filled in) ---
$ws.Range("B7").Value = 1267
$ws.Range("C7").Value = 3802
$ws.Range("B8").Value = 789
$ws.Range("C8").Value = 3981

# --- Column D ("Total de Sacos") formulas ---
# D2 gets its own (non-shared) "B2+C2" formula.
$ws.Range("D2").Formula = "= B2+C2"

# D3:D6 become one shared-formula group using "B3+C3" style refs.
$ws.Range("D3:D6").Formula = "= B3+C3"

# D7:D13 become a second shared-formula group using SUM(B7:C7) style refs.
$ws.Range("D7:D13").Formula = "=SUM(B7:C7)"

# --- Row 14 totals ---
# Drop the old B14/C14 subtotal formulas entirely.
$ws.Range("B14:C14").ClearContents()

# D14 now sums the column D detail cells directly instead of B14:C14.
$ws.Range("D14").Formula = "=SUM(D2,D3,D4,D5,D6,D7,D8,D9,D10,D11,D12,D13)"

# --- Selection bookkeeping (matches the saved sheetView selection in the diff) ---
$ws.Range("J10").Select()
